# Apply a left-rotation of columns B,C,D (new B = old C, new C = old D, new D = old B)
# for rows 6-10 on the "Sheet11" worksheet, leaving column A untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet11")

for ($r = 6; $r -le 10; $r++) {
    $oldB = $ws.Cells.Item($r, 2).Value2
    $oldC = $ws.Cells.Item($r, 3).Value2
    $oldD = $ws.Cells.Item($r, 4).Value2

    $ws.Cells.Item($r, 2).Value2 = $oldC
    $ws.Cells.Item($r, 3).Value2 = $oldD
    $ws.Cells.Item($r, 4).Value2 = $oldB
}
